$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Safe (non-numeric-looking) text values: direct assignment preserves text type.
$ws.Range("D2").Value = "27.337.15"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.832.61"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "1.884.73"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "27.535.97"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "2.091.28"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +0.14%  "

# Numeric-looking text values: Excel would auto-convert these to numbers via a
# direct .Value assignment (matching real Range.Value semantics), which would
# flip the stored cell type from text to number. To preserve the original
# inline-string/text type exactly, build each value as a text-formula result in
# a scratch cell, then copy/paste-special as values into the target cell.
$ws.Range("ZZ1").Formula = "=""1.012"""
$ws.Range("ZZ1").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""314.54"""
$ws.Range("ZZ1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.4745"""
$ws.Range("ZZ1").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.3687"""
$ws.Range("ZZ1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.07459"""
$ws.Range("ZZ1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.8848"""
$ws.Range("ZZ1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.07331"""
$ws.Range("ZZ1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""93.15"""
$ws.Range("ZZ1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""6.583"""
$ws.Range("ZZ1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.000008799"""
$ws.Range("ZZ1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""5.294"""
$ws.Range("ZZ1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""1.892"""
$ws.Range("ZZ1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""151.93"""
$ws.Range("ZZ1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""2.139"""
$ws.Range("ZZ1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""5.248"""
$ws.Range("ZZ1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""117.42"""
$ws.Range("ZZ1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.08996"""
$ws.Range("ZZ1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.7550"""
$ws.Range("ZZ1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""1.176"""
$ws.Range("ZZ1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""4.548"""
$ws.Range("ZZ1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""2.950"""
$ws.Range("ZZ1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""1.103"""
$ws.Range("ZZ1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.05346"""
$ws.Range("ZZ1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.01955"""
$ws.Range("ZZ1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""2.975"""
$ws.Range("ZZ1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""7.288"""
$ws.Range("ZZ1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""2.392"""
$ws.Range("ZZ1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""0.5317"""
$ws.Range("ZZ1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""8.477"""
$ws.Range("ZZ1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""10.56"""
$ws.Range("ZZ1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("ZZ1").Formula = "=""104.96"""
$ws.Range("ZZ1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
